$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 36: add " -BackStore-" header in column B, shift old C/D content left ---
$ws.Range("A36").Copy()
$ws.Range("B36").PasteSpecial(-4122)
$ws.Range("B36").Value = " -BackStore-"
$ws.Range("C36").Value = "Rien pour le moment..."
$ws.Range("D36").ClearContents()

# --- Row 37: add " -Technique-" header in column B (no explicit style), shift old C/D content left ---
$ws.Range("B37").Value = " -Technique-"
$ws.Range("C37").Value = "Rien pour le moment..."
$ws.Range("D37").ClearContents()

# --- Row 38: move " -Encaissement-" header from C to B, shift old C/D content left ---
$ws.Range("A38").Copy()
$ws.Range("B38").PasteSpecial(-4122)
$ws.Range("B38").Value = " -Encaissement-"
$ws.Range("C38").Value = "Rien pour le moment..."
$ws.Range("D38").ClearContents()

$excel.CutCopyMode = $false

# --- Sheet view: scroll position / selection changed ---
$ws.Range("D38").Select()
$ws.Application.ActiveWindow.ScrollRow = 30
